# Update ligand/receptor TPM-derived metrics in the LR-pairs sheet
# (Plau-Plaur) with freshly recomputed values, per "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.44975366666667
$ws.Range("H2").Value = 46.349261
$ws.Range("I2").Value = 0.1287486886000874
$ws.Range("J2").Value = 0.1287486886000874
$ws.Range("M2").Value = 11.128273
$ws.Range("N2").Value = 33.384819
$ws.Range("O2").Value = 0.1975004092010595
$ws.Range("P2").Value = 0.1975004092010595
$ws.Range("Q2").Value = 171.9290765854177
$ws.Range("R2").Value = 1547.361689268759
$ws.Range("S2").Value = 0.02542791868261705
$ws.Range("T2").Value = 0.02542791868261704
$ws.Range("G3").Value = 15.44975366666667
$ws.Range("H3").Value = 46.349261
$ws.Range("I3").Value = 0.1287486886000874
$ws.Range("J3").Value = 0.1287486886000874
$ws.Range("O3").Value = 0.1348505024127659
$ws.Range("P3").Value = 0.1348505024127659
$ws.Range("Q3").Value = 117.3907560530872
$ws.Range("R3").Value = 1056.516804477785
$ws.Range("S3").Value = 0.01736182534270653
$ws.Range("T3").Value = 0.01736182534270653
$ws.Range("G4").Value = 15.44975366666667
$ws.Range("H4").Value = 46.349261
$ws.Range("I4").Value = 0.1287486886000874
$ws.Range("J4").Value = 0.1287486886000874
$ws.Range("M4").Value = 8.642352000000001
$ws.Range("N4").Value = 25.927056
$ws.Range("O4").Value = 0.153381217054937
$ws.Range("P4").Value = 0.153381217054937
$ws.Range("Q4").Value = 133.522209500624
$ws.Range("R4").Value = 1201.699885505616
$ws.Range("S4").Value = 0.0197476305517085
$ws.Range("T4").Value = 0.01974763055170849
$ws.Range("G5").Value = 15.44975366666667
$ws.Range("H5").Value = 46.349261
$ws.Range("I5").Value = 0.1287486886000874
$ws.Range("J5").Value = 0.1287486886000874
$ws.Range("M5").Value = 28.976716
$ws.Range("N5").Value = 86.930148
$ws.Range("O5").Value = 0.5142678713312377
$ws.Range("P5").Value = 0.5142678713312377
$ws.Range("Q5").Value = 447.6831242689586
$ws.Range("R5").Value = 4029.148118420628
$ws.Range("S5").Value = 0.06621131402305536
$ws.Range("T5").Value = 0.06621131402305534
$ws.Range("I6").Value = 0.4074352211478151
$ws.Range("J6").Value = 0.4074352211478151
$ws.Range("M6").Value = 11.128273
$ws.Range("N6").Value = 33.384819
$ws.Range("O6").Value = 0.1975004092010595
$ws.Range("P6").Value = 0.1975004092010595
$ws.Range("Q6").Value = 544.082911460985
$ws.Range("R6").Value = 4896.746203148865
$ws.Range("S6").Value = 0.08046862289961763
$ws.Range("T6").Value = 0.08046862289961763
$ws.Range("I7").Value = 0.4074352211478151
$ws.Range("J7").Value = 0.4074352211478151
$ws.Range("O7").Value = 0.1348505024127659
$ws.Range("P7").Value = 0.1348505024127659
$ws.Range("S7").Value = 0.05494284427243924
$ws.Range("T7").Value = 0.05494284427243924
$ws.Range("I8").Value = 0.4074352211478151
$ws.Range("J8").Value = 0.4074352211478151
$ws.Range("M8").Value = 8.642352000000001
$ws.Range("N8").Value = 25.927056
$ws.Range("O8").Value = 0.153381217054937
$ws.Range("P8").Value = 0.153381217054937
$ws.Range("Q8").Value = 422.54139865464
$ws.Range("R8").Value = 3802.87258789176
$ws.Range("S8").Value = 0.06249291009069928
$ws.Range("T8").Value = 0.06249291009069927
$ws.Range("I9").Value = 0.4074352211478151
$ws.Range("J9").Value = 0.4074352211478151
$ws.Range("M9").Value = 28.976716
$ws.Range("N9").Value = 86.930148
$ws.Range("O9").Value = 0.5142678713312377
$ws.Range("P9").Value = 0.5142678713312377
$ws.Range("Q9").Value = 1416.72800495262
$ws.Range("R9").Value = 12750.55204457358
$ws.Range("S9").Value = 0.209530843885059
$ws.Range("T9").Value = 0.209530843885059
$ws.Range("G10").Value = 24.32144666666666
$ws.Range("H10").Value = 72.96433999999999
$ws.Range("I10").Value = 0.2026798893205849
$ws.Range("J10").Value = 0.2026798893205849
$ws.Range("M10").Value = 11.128273
$ws.Range("N10").Value = 33.384819
$ws.Range("O10").Value = 0.1975004092010595
$ws.Range("P10").Value = 0.1975004092010595
$ws.Range("Q10").Value = 270.6556982616066
$ws.Range("R10").Value = 2435.90128435446
$ws.Range("S10").Value = 0.04002936107764096
$ws.Range("T10").Value = 0.04002936107764096
$ws.Range("G11").Value = 24.32144666666666
$ws.Range("H11").Value = 72.96433999999999
$ws.Range("I11").Value = 0.2026798893205849
$ws.Range("J11").Value = 0.2026798893205849
$ws.Range("O11").Value = 0.1348505024127659
$ws.Range("P11").Value = 0.1348505024127659
$ws.Range("Q11").Value = 184.7999051703222
$ws.Range("R11").Value = 1663.1991465329
$ws.Range("S11").Value = 0.02733148490384465
$ws.Range("T11").Value = 0.02733148490384465
$ws.Range("G12").Value = 24.32144666666666
$ws.Range("H12").Value = 72.96433999999999
$ws.Range("I12").Value = 0.2026798893205849
$ws.Range("J12").Value = 0.2026798893205849
$ws.Range("M12").Value = 8.642352000000001
$ws.Range("N12").Value = 25.927056
$ws.Range("O12").Value = 0.153381217054937
$ws.Range("P12").Value = 0.153381217054937
$ws.Range("Q12").Value = 210.19450324256
$ws.Range("R12").Value = 1891.75052918304
$ws.Range("S12").Value = 0.03108728809655123
$ws.Range("T12").Value = 0.03108728809655123
$ws.Range("G13").Value = 24.32144666666666
$ws.Range("H13").Value = 72.96433999999999
$ws.Range("I13").Value = 0.2026798893205849
$ws.Range("J13").Value = 0.2026798893205849
$ws.Range("M13").Value = 28.976716
$ws.Range("N13").Value = 86.930148
$ws.Range("O13").Value = 0.5142678713312377
$ws.Range("P13").Value = 0.5142678713312377
$ws.Range("Q13").Value = 704.7556527691465
$ws.Range("R13").Value = 6342.80087492232
$ws.Range("S13").Value = 0.104231755242548
$ws.Range("T13").Value = 0.104231755242548
$ws.Range("G14").Value = 31.33616366666666
$ws.Range("H14").Value = 94.00849099999999
$ws.Range("I14").Value = 0.2611362009315126
$ws.Range("J14").Value = 0.2611362009315126
$ws.Range("M14").Value = 11.128273
$ws.Range("N14").Value = 33.384819
$ws.Range("O14").Value = 0.1975004092010595
$ws.Range("P14").Value = 0.1975004092010595
$ws.Range("Q14").Value = 348.7173840553476
$ws.Range("R14").Value = 3138.456456498129
$ws.Range("S14").Value = 0.05157450654118382
$ws.Range("T14").Value = 0.05157450654118382
$ws.Range("G15").Value = 31.33616366666666
$ws.Range("H15").Value = 94.00849099999999
$ws.Range("I15").Value = 0.2611362009315126
$ws.Range("J15").Value = 0.2611362009315126
$ws.Range("O15").Value = 0.1348505024127659
$ws.Range("P15").Value = 0.1348505024127659
$ws.Range("Q15").Value = 238.0993266300372
$ws.Range("R15").Value = 2142.893939670335
$ws.Range("S15").Value = 0.03521434789377546
$ws.Range("T15").Value = 0.03521434789377546
$ws.Range("G16").Value = 31.33616366666666
$ws.Range("H16").Value = 94.00849099999999
$ws.Range("I16").Value = 0.2611362009315126
$ws.Range("J16").Value = 0.2611362009315126
$ws.Range("M16").Value = 8.642352000000001
$ws.Range("N16").Value = 25.927056
$ws.Range("O16").Value = 0.153381217054937
$ws.Range("P16").Value = 0.153381217054937
$ws.Range("Q16").Value = 270.818156736944
$ws.Range("R16").Value = 2437.363410632496
$ws.Range("S16").Value = 0.04005338831597797
$ws.Range("T16").Value = 0.04005338831597796
$ws.Range("G17").Value = 31.33616366666666
$ws.Range("H17").Value = 94.00849099999999
$ws.Range("I17").Value = 0.2611362009315126
$ws.Range("J17").Value = 0.2611362009315126
$ws.Range("M17").Value = 28.976716
$ws.Range("N17").Value = 86.930148
$ws.Range("O17").Value = 0.5142678713312377
$ws.Range("P17").Value = 0.5142678713312377
$ws.Range("Q17").Value = 908.0191150985186
$ws.Range("R17").Value = 8172.172035886668
$ws.Range("S17").Value = 0.1342939581805754
$ws.Range("T17").Value = 0.1342939581805754
